# Revert "Added converter option":
# The converter option query string used on header cells changes back
# from "?readAs=text" to "?type=string" on the two test sheets that used
# it, and the active-sheet/selection state in the workbook reverts to an
# earlier snapshot (indexedListAsLeafTestLeft becomes the active/selected
# sheet instead of indexedListAsLeaf).

$wb = $excel.ActiveWorkbook

# --- indexedListAsLeafTestOption ---------------------------------------
$wsOption = $wb.Worksheets.Item("indexedListAsLeafTestOption")
$wsOption.Range("C1").Value = "listAsLeafTestOption#test?type=string"
$wsOption.Range("E1").Value = "listAsLeafTestOption#list[1]?type=string"

# --- indexedListAsLeafTestLeft ------------------------------------------
$wsLeft = $wb.Worksheets.Item("indexedListAsLeafTestLeft")
$wsLeft.Range("C1").Value = "listAsLeafTestLeft[0]#test?type=string"
$wsLeft.Range("E1").Value = "listAsLeafTestLeft[0]#list[1]?type=string"
$wsLeft.Range("J1").Value = "listAsLeafTestLeft[1]#list[1]?type=string"

# --- Active sheet / selection state -------------------------------------
# Before the revert, "indexedListAsLeaf" was the active tab with cell E2
# selected; after the revert "indexedListAsLeafTestLeft" is active with H3
# selected. Activating the sheet and selecting the cell updates
# tabSelected / activeTab and the per-sheet selection in one step.
$wsLeft.Activate()
$wsLeft.Range("H3").Select()
